# Update the Metadata sheet:
#  - "Experimental" row (row 7): set the Value column (B7) to "true"
#  - "Date" row (row 8): bump the Value column (B8) to the new timestamp

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B7").Value = "true"
$ws.Range("B8").Value = "2025-07-21T12:46:15+00:00"
